$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing row 12: B12 text and C12 date value
$ws.Range("B12").Value = "Postman5PT2"
$ws.Range("C12").Value = 44958.59311342592

# Add new row 15 with data, copying formatting from row 14 (the last existing row)
# so the styled columns (A = bold/bordered id column, C = date column) keep the
# same cell styles as the rest of the table.
$ws.Range("A14").Copy($ws.Range("A15"))
$ws.Range("C14").Copy($ws.Range("C15"))

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Saving new email 2"
$ws.Range("C15").Value = 44958.59305555555
$ws.Range("D15").Value = "Hope this works AGAIN!"
$ws.Range("E15").Value = "duanevaughn@hotmail.com"
$ws.Range("F15").Value = "Duane Vaughn"
$ws.Range("G15").Value = "{No Recipient Email}"
$ws.Range("H15").Value = "{No Recipient Name}"
$ws.Range("I15").Value = $true
